# Seminar Registration - ATDD Scenarios
# "Added comment to each test codeunit to explain it's scope and updated ATDD sheet"
#
# On the "ATDD Scenarios - Posting only" sheet, mark every scenario /
# given-when-then row (rows 32-94) with an "X" in column D ("Positve-Negative"
# helper column), then leave the selection on the first newly touched block
# (D33:D39) as the author's workbook view shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATDD Scenarios - Posting only")

for ($row = 32; $row -le 94; $row++) {
    $ws.Cells.Item($row, 4).Value = "X"
}

$ws.Activate()
$ws.Range("D33:D39").Select()
